# This script applies the "Make changes to 'The work everyone did'" edit:
#  - Georgi Trendafilov / Readme paragraphs merge into one, with the name
#    split into spell-checked runs and the description rewritten.
#  - Petar Nikolov paragraph: name split into spell-checked runs, new
#    description.
#  - Mario Boshev / "and visual changes" paragraphs merge into one, with
#    the name split into spell-checked runs and a new description.
#  - Georgi Ivanov paragraph: trailing runs consolidated into one run
#    (no visible text change).

$d = $word.ActiveDocument
$ndash = [char]8211
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Find-ParaIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) "Georgi Trendafilov - Documentation, Presentation and" + "Readme"
#    -> merge into one paragraph:
#    "Georgi Trendafilov - Readme and the HTML structure of the site"
# ---------------------------------------------------------------------------
$idx1 = Find-ParaIndex("Georgi Trendafilov $ndash Documentation")
$p1a = $d.Paragraphs.Item($idx1)
$p1b = $d.Paragraphs.Item($idx1 + 1)
$range1 = $d.Range($p1a.Range.Start, $p1b.Range.End)

$xml1 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val='NoSpacing'/>
    <w:numPr>
      <w:ilvl w:val='0'/>
      <w:numId w:val='12'/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val='EA4E4E' w:themeColor='accent1'/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t xml:space='preserve'>Georgi </w:t>
  </w:r>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val='EA4E4E' w:themeColor='accent1'/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t>Trendafilov</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val='EA4E4E' w:themeColor='accent1'/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t xml:space='preserve'> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t xml:space='preserve'>$ndash </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t>Readme and the HTML structure of the site</w:t>
  </w:r>
</w:p>
"@

$range1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) "Petar Nikolov - Main body and construction of the website"
#    -> "Petar Nikolov - Documentation of the project, finding information
#        about the crew equipment"
# ---------------------------------------------------------------------------
$idx2 = Find-ParaIndex("Petar Nikolov $ndash Main body")
$p2 = $d.Paragraphs.Item($idx2)

$xml2 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val='NoSpacing'/>
    <w:numPr>
      <w:ilvl w:val='0'/>
      <w:numId w:val='12'/>
    </w:numPr>
    <w:spacing w:line='276' w:lineRule='auto'/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val='EA4E4E' w:themeColor='accent1'/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t>Petar</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val='EA4E4E' w:themeColor='accent1'/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t xml:space='preserve'> </w:t>
  </w:r>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val='EA4E4E' w:themeColor='accent1'/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t>Nikolov</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t xml:space='preserve'> $ndash </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t>Documentation of the project, finding information about the crew equipment</w:t>
  </w:r>
</w:p>
"@

$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) "Mario Boshev - Improvements to the site, photos " + "and visual changes"
#    -> merge into one paragraph:
#    "Mario Boshev - Presentation of the project, finding most of the images"
# ---------------------------------------------------------------------------
$idx3 = Find-ParaIndex("Mario Boshev $ndash Improvements")
$p3a = $d.Paragraphs.Item($idx3)
$p3b = $d.Paragraphs.Item($idx3 + 1)
$range3 = $d.Range($p3a.Range.Start, $p3b.Range.End)

$xml3 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val='NoSpacing'/>
    <w:numPr>
      <w:ilvl w:val='0'/>
      <w:numId w:val='12'/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val='EA4E4E' w:themeColor='accent1'/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t xml:space='preserve'>Mario </w:t>
  </w:r>
  <w:proofErr w:type='spellStart'/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val='EA4E4E' w:themeColor='accent1'/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t>Boshev</w:t>
  </w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val='EA4E4E' w:themeColor='accent1'/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t xml:space='preserve'> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t xml:space='preserve'>$ndash </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t>Presentation of the project, finding most of the images</w:t>
  </w:r>
</w:p>
"@

$range3.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 4) "Georgi Ivanov - Most of the CSS code, arrangement of the "
#    -> same text, trailing runs consolidated into a single run.
# ---------------------------------------------------------------------------
$idx4 = Find-ParaIndex("Georgi Ivanov $ndash Most of the CSS")
$p4 = $d.Paragraphs.Item($idx4)

$xml4 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val='NoSpacing'/>
    <w:numPr>
      <w:ilvl w:val='0'/>
      <w:numId w:val='12'/>
    </w:numPr>
    <w:rPr>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val='EA4E4E' w:themeColor='accent1'/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t xml:space='preserve'>Georgi Ivanov </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val='32'/>
      <w:szCs w:val='32'/>
    </w:rPr>
    <w:t xml:space='preserve'>$ndash Most of the CSS code, arrangement of the </w:t>
  </w:r>
</w:p>
"@

$p4.Range.InsertXML($xml4)

Write-Host "Edits applied successfully."
